$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26
$ws.Cells.Item(26, 1).Value = 112076816
$ws.Cells.Item(26, 2).Value = 96735
$ws.Cells.Item(26, 4).Value = 'VU'
$ws.Cells.Item(26, 5).Value = 220787
$ws.Cells.Item(26, 6).Value = 'Knärot'
$ws.Cells.Item(26, 7).Value = 'Goodyera repens'
$ws.Cells.Item(26, 8).Value = '(L.) R. Br.'
$ws.Cells.Item(26, 17).Value = 485618
$ws.Cells.Item(26, 18).Value = 7005614

# Row 27
$ws.Cells.Item(27, 1).Value = 112076820
$ws.Cells.Item(27, 2).Value = 98980
$ws.Cells.Item(27, 4).Value = 'LC'
$ws.Cells.Item(27, 5).Value = 222498
$ws.Cells.Item(27, 6).Value = 'Blåsippa'
$ws.Cells.Item(27, 7).Value = 'Hepatica nobilis'
$ws.Cells.Item(27, 8).Value = 'Schreb.'
$ws.Cells.Item(27, 17).Value = 485536
$ws.Cells.Item(27, 18).Value = 7005851

# Row 28
$ws.Cells.Item(28, 1).Value = 112076813
$ws.Cells.Item(28, 2).Value = 89553
$ws.Cells.Item(28, 4).Value = 'NT'
$ws.Cells.Item(28, 5).Value = 1202
$ws.Cells.Item(28, 6).Value = 'Ullticka'
$ws.Cells.Item(28, 7).Value = 'Phellinidium ferrugineofuscum'
$ws.Cells.Item(28, 8).Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Cells.Item(28, 17).Value = 485752
$ws.Cells.Item(28, 18).Value = 7005707

# Row 29
$ws.Cells.Item(29, 1).Value = 112076817
$ws.Cells.Item(29, 2).Value = 96735
$ws.Cells.Item(29, 4).Value = 'VU'
$ws.Cells.Item(29, 5).Value = 220787
$ws.Cells.Item(29, 6).Value = 'Knärot'
$ws.Cells.Item(29, 7).Value = 'Goodyera repens'
$ws.Cells.Item(29, 8).Value = '(L.) R. Br.'
$ws.Cells.Item(29, 17).Value = 485596
$ws.Cells.Item(29, 18).Value = 7005613

# Row 30
$ws.Cells.Item(30, 1).Value = 112076818
$ws.Cells.Item(30, 2).Value = 89047
$ws.Cells.Item(30, 4).Value = 'NT'
$ws.Cells.Item(30, 5).Value = 3286
$ws.Cells.Item(30, 6).Value = 'Flattoppad klubbsvamp'
$ws.Cells.Item(30, 7).Value = 'Clavariadelphus truncatus'
$ws.Cells.Item(30, 8).Value = '(Quél.) Donk'
$ws.Cells.Item(30, 17).Value = 485662
$ws.Cells.Item(30, 18).Value = 7005637

# Row 31
$ws.Cells.Item(31, 1).Value = 112076812
$ws.Cells.Item(31, 2).Value = 86371
$ws.Cells.Item(31, 4).Value = 'NT'
$ws.Cells.Item(31, 5).Value = 4412
$ws.Cells.Item(31, 6).Value = 'Äggvaxskivling'
$ws.Cells.Item(31, 7).Value = 'Hygrophorus karstenii'
$ws.Cells.Item(31, 8).Value = 'Sacc. & Cub.'
$ws.Cells.Item(31, 17).Value = 485781
$ws.Cells.Item(31, 18).Value = 7005721

# Row 32
$ws.Cells.Item(32, 1).Value = 112076814
$ws.Cells.Item(32, 2).Value = 90235
$ws.Cells.Item(32, 4).Value = 'LC'
$ws.Cells.Item(32, 5).Value = 3298
$ws.Cells.Item(32, 6).Value = 'Trådticka'
$ws.Cells.Item(32, 7).Value = 'Climacocystis borealis'
$ws.Cells.Item(32, 8).Value = '(Fr.) Kotl. & Pouzar'
$ws.Cells.Item(32, 17).Value = 485714
$ws.Cells.Item(32, 18).Value = 7005798

# Row 33
$ws.Cells.Item(33, 1).Value = 112076819
$ws.Cells.Item(33, 2).Value = 98980
$ws.Cells.Item(33, 4).Value = 'LC'
$ws.Cells.Item(33, 5).Value = 222498
$ws.Cells.Item(33, 6).Value = 'Blåsippa'
$ws.Cells.Item(33, 7).Value = 'Hepatica nobilis'
$ws.Cells.Item(33, 8).Value = 'Schreb.'
$ws.Cells.Item(33, 17).Value = 485664
$ws.Cells.Item(33, 18).Value = 7005761

# Row 34
$ws.Cells.Item(34, 1).Value = 112076815
$ws.Cells.Item(34, 2).Value = 96735
$ws.Cells.Item(34, 4).Value = 'VU'
$ws.Cells.Item(34, 5).Value = 220787
$ws.Cells.Item(34, 6).Value = 'Knärot'
$ws.Cells.Item(34, 7).Value = 'Goodyera repens'
$ws.Cells.Item(34, 8).Value = '(L.) R. Br.'
$ws.Cells.Item(34, 17).Value = 485636
$ws.Cells.Item(34, 18).Value = 7005629

# Row 35
$ws.Cells.Item(35, 1).Value = 112076811
$ws.Cells.Item(35, 2).Value = 86371
$ws.Cells.Item(35, 4).Value = 'NT'
$ws.Cells.Item(35, 5).Value = 4412
$ws.Cells.Item(35, 6).Value = 'Äggvaxskivling'
$ws.Cells.Item(35, 7).Value = 'Hygrophorus karstenii'
$ws.Cells.Item(35, 8).Value = 'Sacc. & Cub.'
$ws.Cells.Item(35, 17).Value = 485716
$ws.Cells.Item(35, 18).Value = 7005807
